# Saved Booma's criteria for keeping potential matches
#
# Row 8 previously carried example-compare text in columns M and N
# ("John Doe, 10-10-2019, 123 Main") that didn't actually apply to that
# rule (drop on name/dob mismatch + single dob element match). Clear
# those two stray example cells out of the worksheet entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M8:N8").Clear()

# Reflect the updated selection/zoom the author ended the session with.
$excel.ActiveWindow.Zoom = 140
$ws.Range("P8").Select()
